$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2017-02-17 07:51:09"

$wsZhCn.Range("H4").Value = "2017-02-17 07:50:51"
$wsZhCn.Range("L4").Value = "2017-02-17 07:51:47"

$wsDeDe.Range("H4").Value = "2017-02-17 07:51:09"
$wsDeDe.Range("L4").Value = "2017-02-17 07:52:10"
